$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4444.9
$ws.Range("I33").Value = 275.16666
$ws.Range("J33").Value = 10699.5
$ws.Range("K33").Value = 275.16666
$ws.Range("L33").Value = 10699.5
$ws.Range("M33").Value = -46.16665999999998
$ws.Range("N33").Value = -11157.5
$ws.Range("H39").Value = 314.33334
$ws.Range("I39").Value = 314.33334
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 943.0000200000001
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -647.0000200000001
$ws.Range("N39").ClearContents()
$ws.Range("H76").Value = 9220.444
$ws.Range("I76").Value = 7985
$ws.Range("J76").Value = 9374.875
$ws.Range("K76").Value = 7985
$ws.Range("L76").Value = 9374.875
$ws.Range("M76").Value = -7670
$ws.Range("N76").Value = -10004.875
$ws.Range("H79").Value = 9220.444
$ws.Range("I79").Value = 7985
$ws.Range("J79").Value = 9374.875
$ws.Range("K79").Value = 7985
$ws.Range("L79").Value = 9374.875
$ws.Range("M79").Value = -6893
$ws.Range("N79").Value = -11558.875
$ws.Range("H86").Value = 5618.75
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 5825
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 5825
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -8071
$ws.Range("H89").Value = 5618.75
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 5825
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 29125
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -40357
$ws.Range("H98").Value = 11137.2
$ws.Range("J98").Value = 50749.5
$ws.Range("L98").Value = 50749.5
$ws.Range("N98").Value = -53745.5
$ws.Range("H103").Value = 5228.2856
$ws.Range("I103").Value = 899.5
$ws.Range("J103").Value = 6959.8
$ws.Range("K103").Value = 2698.5
$ws.Range("L103").Value = 20879.4
$ws.Range("M103").Value = -2112.5
$ws.Range("N103").Value = -22051.4
$ws.Range("H116").Value = 8083.3335
$ws.Range("J116").Value = 8901
$ws.Range("L116").Value = 8901
$ws.Range("N116").Value = -15785
$ws.Range("H122").Value = 11137.2
$ws.Range("J122").Value = 50749.5
$ws.Range("L122").Value = 152248.5
$ws.Range("N122").Value = -157148.5
$ws.Range("H127").Value = 863.3333
$ws.Range("J127").Value = 4499
$ws.Range("L127").Value = 13497
$ws.Range("N127").Value = -23417
$ws.Range("H137").Value = 4134.8438
$ws.Range("I137").Value = 3316.75
$ws.Range("J137").Value = 5498.3335
$ws.Range("K137").Value = 9950.25
$ws.Range("L137").Value = 16495.0005
$ws.Range("M137").Value = -7400.25
$ws.Range("N137").Value = -21595.0005
$ws.Range("H138").Value = 2951.9768
$ws.Range("I138").Value = 2479.3333
$ws.Range("J138").Value = 3403.1365
$ws.Range("K138").Value = 7437.999899999999
$ws.Range("L138").Value = 10209.4095
$ws.Range("M138").Value = -2297.999899999999
$ws.Range("N138").Value = -20489.4095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6515.8
$ws.Range("I61").Value = 4950.05
$ws.Range("K61").Value = 4950.05
$ws.Range("M61").Value = -4738.05
$ws.Range("H74").Value = 17549938
$ws.Range("I74").Value = 25646548
$ws.Range("K74").Value = 25646548
$ws.Range("M74").Value = -25645674
$ws.Range("H77").Value = 17549938
$ws.Range("I77").Value = 25646548
$ws.Range("K77").Value = 128232740
$ws.Range("M77").Value = -128228372
$ws.Range("H97").Value = 8933.333000000001
$ws.Range("I97").Value = 8912.25
$ws.Range("J97").Value = 8975.5
$ws.Range("K97").Value = 8912.25
$ws.Range("L97").Value = 8975.5
$ws.Range("M97").Value = -8416.25
$ws.Range("N97").Value = -9967.5
$ws.Range("H112").Value = 68841.336
$ws.Range("J112").Value = 68841.336
$ws.Range("L112").Value = 68841.336
$ws.Range("N112").Value = -71795.336
$ws.Range("H122").Value = 4286.8096
$ws.Range("I122").Value = 3668.2
$ws.Range("K122").Value = 11004.6
$ws.Range("M122").Value = -8554.599999999999
$ws.Range("H132").Value = 4847.22
$ws.Range("I132").Value = 3530.282
$ws.Range("K132").Value = 10590.846
$ws.Range("M132").Value = -8060.846000000001
$ws.Range("H136").Value = 6515.8
$ws.Range("I136").Value = 4950.05
$ws.Range("K136").Value = 14850.15
$ws.Range("M136").Value = -12300.15

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 18400.666
$ws.Range("I105").Value = 19003.416
$ws.Range("K105").Value = 19003.416
$ws.Range("M105").Value = -17256.416
$ws.Range("H107").Value = 3600.2083
$ws.Range("I107").Value = 3411.25
$ws.Range("J107").Value = 3978.125
$ws.Range("K107").Value = 3411.25
$ws.Range("L107").Value = 3978.125
$ws.Range("M107").Value = -1491.25
$ws.Range("N107").Value = -7818.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6999.6665
$ws.Range("J16").Value = 6999.6665
$ws.Range("L16").Value = 6999.6665
$ws.Range("N16").Value = -7573.6665
$ws.Range("H105").Value = 3018.5
$ws.Range("J105").Value = 2896.2144
$ws.Range("L105").Value = 2896.2144
$ws.Range("N105").Value = -6390.2144
$ws.Range("H113").Value = 6999.6665
$ws.Range("J113").Value = 6999.6665
$ws.Range("L113").Value = 6999.6665
$ws.Range("N113").Value = -11339.6665
$ws.Range("H122").Value = 10612.2
$ws.Range("I122").Value = 3335.3333
$ws.Range("J122").Value = 21527.5
$ws.Range("K122").Value = 10005.9999
$ws.Range("L122").Value = 64582.5
$ws.Range("M122").Value = -7555.999899999999
$ws.Range("N122").Value = -69482.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2736348.2
$ws.Range("I4").Value = 3608277.5
$ws.Range("K4").Value = 10824832.5
$ws.Range("M4").Value = -10824720.5
$ws.Range("H55").Value = 1344.8
$ws.Range("I55").Value = 287
$ws.Range("J55").Value = 2050
$ws.Range("K55").Value = 861
$ws.Range("L55").Value = 6150
$ws.Range("M55").Value = -684
$ws.Range("N55").Value = -6504
$ws.Range("H64").Value = 34487156
$ws.Range("J64").Value = 5214.476
$ws.Range("L64").Value = 15643.428
$ws.Range("N64").Value = -16183.428
$ws.Range("H67").Value = 34487156
$ws.Range("J67").Value = 5214.476
$ws.Range("L67").Value = 15643.428
$ws.Range("N67").Value = -17515.428
$ws.Range("H97").Value = 978.8333
$ws.Range("I97").Value = 250
$ws.Range("J97").Value = 1343.25
$ws.Range("K97").Value = 750
$ws.Range("L97").Value = 4029.75
$ws.Range("M97").Value = -254
$ws.Range("N97").Value = -5021.75
$ws.Range("H137").Value = 1001605.7
$ws.Range("I137").Value = 1429167
$ws.Range("J137").Value = 3962.6667
$ws.Range("K137").Value = 4287501
$ws.Range("L137").Value = 11888.0001
$ws.Range("M137").Value = -4282401
$ws.Range("N137").Value = -22088.0001
$ws.Range("H140").Value = 3376.524
$ws.Range("I140").Value = 2826.7222
$ws.Range("K140").Value = 8480.1666
$ws.Range("M140").Value = -3300.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 4217.2
$ws.Range("I132").Value = 3391.1177
$ws.Range("J132").Value = 5972.625
$ws.Range("K132").Value = 10173.3531
$ws.Range("L132").Value = 17917.875
$ws.Range("M132").Value = -7643.3531
$ws.Range("N132").Value = -22977.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15852.667
$ws.Range("I7").Value = 7657
$ws.Range("J7").Value = 44537.5
$ws.Range("K7").Value = 7657
$ws.Range("L7").Value = 44537.5
$ws.Range("M7").Value = -7545
$ws.Range("N7").Value = -44761.5
$ws.Range("H22").Value = 29407.846
$ws.Range("I22").Value = 27475.312
$ws.Range("K22").Value = 27475.312
$ws.Range("M22").Value = -27180.312
$ws.Range("H27").Value = 29407.846
$ws.Range("I27").Value = 27475.312
$ws.Range("K27").Value = 27475.312
$ws.Range("M27").Value = -27368.312
$ws.Range("H93").Value = 14288.723
$ws.Range("I93").Value = 11775.228
$ws.Range("J93").Value = 18238.5
$ws.Range("K93").Value = 11775.228
$ws.Range("L93").Value = 18238.5
$ws.Range("M93").Value = -10527.228
$ws.Range("N93").Value = -20734.5
$ws.Range("H126").Value = 15852.667
$ws.Range("I126").Value = 7657
$ws.Range("J126").Value = 44537.5
$ws.Range("K126").Value = 22971
$ws.Range("L126").Value = 133612.5
$ws.Range("M126").Value = -20501
$ws.Range("N126").Value = -138552.5
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
$ws.Range("H136").Value = 6269.643
$ws.Range("I136").Value = 3940.238
$ws.Range("K136").Value = 11820.714
$ws.Range("M136").Value = -9270.714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4028.3333
$ws.Range("J126").Value = 4718.8125
$ws.Range("L126").Value = 14156.4375
$ws.Range("N126").Value = -19096.4375
$ws.Range("H136").Value = 3084.0857
$ws.Range("I136").Value = 1441.7368
$ws.Range("K136").Value = 4325.2104
$ws.Range("M136").Value = -1775.2104
